$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.614.86"
$ws.Range("E2").Value = "  +4.76%  "

$ws.Range("D3").Value = "1.918.03"
$ws.Range("E3").Value = "  +3.28%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'334.56"
$ws.Range("E5").Value = "  +1.73%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.4667"
$ws.Range("E7").Value = "  +2.41%  "

$ws.Range("D8").Value = "'0.4114"
$ws.Range("E8").Value = "  +4.71%  "

$ws.Range("D9").Value = "'48.12"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").Value = "'0.08034"
$ws.Range("E10").Value = "  +2.64%  "

$ws.Range("D11").Value = "'1.013"
$ws.Range("E11").Value = "  +3.32%  "

$ws.Range("D12").Value = "'22.34"
$ws.Range("E12").Value = "  +5.13%  "

$ws.Range("D13").Value = "1.898.73"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").Value = "'5.993"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").Value = "'7.173"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").Value = "'90.08"

$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "'0.06597"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("D20").Value = "'17.84"
$ws.Range("E20").Value = "  +4.60%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "29.606.45"
$ws.Range("E22").Value = "  +4.76%  "

$ws.Range("D23").Value = "'5.581"
$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("D24").Value = "'11.66"
$ws.Range("E24").Value = "  +7.74%  "

$ws.Range("D25").Value = "'2.215"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("D26").Value = "2.195.60"
$ws.Range("E26").Value = "  +5.96%  "

$ws.Range("D27").Value = "'156.08"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").Value = "'19.89"
$ws.Range("E28").Value = "  +3.65%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.741"
$ws.Range("E29").Value = "  +8.42%  "

$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'2.136"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").Value = "'117.51"
$ws.Range("E31").Value = "  +0.97%  "

$ws.Range("E32").Value = "  +13.43%  "

$ws.Range("D33").Value = "'0.09475"
$ws.Range("E33").Value = "  +2.18%  "

$ws.Range("D34").Value = "'1.431"
$ws.Range("E34").Value = "  +3.95%  "

$ws.Range("D35").Value = "'3.580"
$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").Value = "'5.402"
$ws.Range("E36").Value = "  +3.87%  "

$ws.Range("D37").Value = "'0.06133"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("D38").Value = "'0.02269"
$ws.Range("E38").Value = "  +2.93%  "

$ws.Range("D39").Value = "'8.425"
$ws.Range("E39").Value = "  +2.20%  "

$ws.Range("D40").Value = "'1.183"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("D41").Value = "'0.5902"
$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1845"
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'10.22"
$ws.Range("E44").Value = "  +2.43%  "

$ws.Range("D45").Value = "'1.262"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'2.362"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.07516"
$ws.Range("E47").Value = "  +5.03%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'12.25"
$ws.Range("E48").Value = "  +3.69%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.5571"
$ws.Range("E49").Value = "  +3.17%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.933"
$ws.Range("E50").Value = "  +3.57%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'112.88"
$ws.Range("E51").Value = "  +2.69%  "
